# Update iNEXT-derived diversity estimates (Table_S4) after package updates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.1716685820065127
$ws.Range("C2").Value = 2.00124302412812
$ws.Range("D2").Value = 4.718067388696906
$ws.Range("E2").Value = 1.351231779597265
$ws.Range("F2").Value = 5.398287328394347
$ws.Range("G2").Value = 1.2872725023091371

$ws.Range("B3").Value = 1.5284257587900505
$ws.Range("C3").Value = 0.9849155338944934
$ws.Range("D3").Value = 3.5999748924387234
$ws.Range("E3").Value = 0.9203589835436573
$ws.Range("F3").Value = 4.247972640506182
$ws.Range("G3").Value = 0.9774366686076593

$ws.Range("B4").Value = 5.878683605029869
$ws.Range("C4").Value = 0.8023107209044661
$ws.Range("D4").Value = 5.008172886841643
$ws.Range("E4").Value = 0.4578169882560387
$ws.Range("F4").Value = 4.396447212636887
$ws.Range("G4").Value = 0.7869921693224223

$ws.Range("B5").Value = 4.65059330644514
$ws.Range("C5").Value = 0.7913677463237762
$ws.Range("D5").Value = 4.370886260119848
$ws.Range("E5").Value = 1.061707234163964
$ws.Range("F5").Value = 4.136227690441243
$ws.Range("G5").Value = 1.0002426955902959
